$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update Maximo (C2) ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 693.987917612715

# --- Sheet "Solucion": reshuffle Pedido (A) and Salida (B) columns ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsSolucion.Range("A2").Value = "Pedido_17"
$wsSolucion.Range("B2").Value = "S001"
$wsSolucion.Range("A3").Value = "Pedido_43"
$wsSolucion.Range("B3").Value = "S025"
$wsSolucion.Range("A4").Value = "Pedido_74"
$wsSolucion.Range("B4").Value = "S041"
$wsSolucion.Range("A5").Value = "Pedido_75"
$wsSolucion.Range("B5").Value = "S065"
$wsSolucion.Range("A6").Value = "Pedido_59"
$wsSolucion.Range("B6").Value = "S069"
$wsSolucion.Range("A7").Value = "Pedido_8"
$wsSolucion.Range("B7").Value = "S005"
$wsSolucion.Range("A8").Value = "Pedido_3"
$wsSolucion.Range("B8").Value = "S029"
$wsSolucion.Range("A9").Value = "Pedido_37"
$wsSolucion.Range("B9").Value = "S045"
$wsSolucion.Range("A10").Value = "Pedido_78"
$wsSolucion.Range("B10").Value = "S002"
$wsSolucion.Range("A11").Value = "Pedido_29"
$wsSolucion.Range("B11").Value = "S066"
$wsSolucion.Range("A12").Value = "Pedido_16"
$wsSolucion.Range("B12").Value = "S042"
$wsSolucion.Range("A13").Value = "Pedido_80"
$wsSolucion.Range("B13").Value = "S006"
$wsSolucion.Range("A14").Value = "Pedido_67"
$wsSolucion.Range("B14").Value = "S026"
$wsSolucion.Range("A15").Value = "Pedido_30"
$wsSolucion.Range("B15").Value = "S046"
$wsSolucion.Range("A16").Value = "Pedido_52"
$wsSolucion.Range("B16").Value = "S070"
$wsSolucion.Range("A17").Value = "Pedido_79"
$wsSolucion.Range("B17").Value = "S003"
$wsSolucion.Range("A18").Value = "Pedido_34"
$wsSolucion.Range("B18").Value = "S043"
$wsSolucion.Range("A19").Value = "Pedido_4"
$wsSolucion.Range("B19").Value = "S030"
$wsSolucion.Range("A20").Value = "Pedido_45"
$wsSolucion.Range("B20").Value = "S007"
$wsSolucion.Range("A21").Value = "Pedido_33"
$wsSolucion.Range("B21").Value = "S067"
$wsSolucion.Range("A22").Value = "Pedido_40"
$wsSolucion.Range("B22").Value = "S027"
$wsSolucion.Range("A23").Value = "Pedido_18"
$wsSolucion.Range("B23").Value = "S047"
$wsSolucion.Range("A24").Value = "Pedido_21"
$wsSolucion.Range("B24").Value = "S004"
$wsSolucion.Range("A25").Value = "Pedido_46"
$wsSolucion.Range("B25").Value = "S031"
$wsSolucion.Range("A26").Value = "Pedido_5"
$wsSolucion.Range("B26").Value = "S071"
$wsSolucion.Range("A27").Value = "Pedido_32"
$wsSolucion.Range("B27").Value = "S044"
$wsSolucion.Range("A28").Value = "Pedido_69"
$wsSolucion.Range("B28").Value = "S028"
$wsSolucion.Range("A29").Value = "Pedido_50"
$wsSolucion.Range("B29").Value = "S008"
$wsSolucion.Range("A30").Value = "Pedido_64"
$wsSolucion.Range("B30").Value = "S032"
$wsSolucion.Range("A31").Value = "Pedido_62"
$wsSolucion.Range("B31").Value = "S048"
$wsSolucion.Range("A32").Value = "Pedido_47"
$wsSolucion.Range("B32").Value = "S068"
$wsSolucion.Range("A33").Value = "Pedido_66"
$wsSolucion.Range("B33").Value = "S009"
$wsSolucion.Range("A34").Value = "Pedido_10"
$wsSolucion.Range("B34").Value = "S049"
$wsSolucion.Range("A35").Value = "Pedido_22"
$wsSolucion.Range("B35").Value = "S013"
$wsSolucion.Range("A36").Value = "Pedido_49"
$wsSolucion.Range("B36").Value = "S072"
$wsSolucion.Range("A37").Value = "Pedido_42"
$wsSolucion.Range("B37").Value = "S033"
$wsSolucion.Range("A38").Value = "Pedido_63"
$wsSolucion.Range("B38").Value = "S010"
$wsSolucion.Range("A39").Value = "Pedido_35"
$wsSolucion.Range("B39").Value = "S037"
$wsSolucion.Range("A40").Value = "Pedido_39"
$wsSolucion.Range("B40").Value = "S053"
$wsSolucion.Range("A41").Value = "Pedido_13"
$wsSolucion.Range("B41").Value = "S073"
$wsSolucion.Range("A42").Value = "Pedido_77"
$wsSolucion.Range("B42").Value = "S014"
$wsSolucion.Range("A43").Value = "Pedido_25"
$wsSolucion.Range("B43").Value = "S034"
$wsSolucion.Range("A44").Value = "Pedido_57"
$wsSolucion.Range("B44").Value = "S050"
$wsSolucion.Range("A45").Value = "Pedido_73"
$wsSolucion.Range("B45").Value = "S011"
$wsSolucion.Range("A46").Value = "Pedido_71"
$wsSolucion.Range("B46").Value = "S077"
$wsSolucion.Range("A47").Value = "Pedido_65"
$wsSolucion.Range("B47").Value = "S054"
$wsSolucion.Range("A48").Value = "Pedido_12"
$wsSolucion.Range("B48").Value = "S038"
$wsSolucion.Range("A49").Value = "Pedido_36"
$wsSolucion.Range("B49").Value = "S015"
$wsSolucion.Range("A50").Value = "Pedido_1"
$wsSolucion.Range("B50").Value = "S051"
$wsSolucion.Range("A51").Value = "Pedido_23"
$wsSolucion.Range("B51").Value = "S074"
$wsSolucion.Range("A52").Value = "Pedido_58"
$wsSolucion.Range("B52").Value = "S055"
$wsSolucion.Range("A53").Value = "Pedido_61"
$wsSolucion.Range("B53").Value = "S012"
$wsSolucion.Range("A54").Value = "Pedido_28"
$wsSolucion.Range("B54").Value = "S035"
$wsSolucion.Range("A55").Value = "Pedido_51"
$wsSolucion.Range("B55").Value = "S078"
$wsSolucion.Range("A56").Value = "Pedido_48"
$wsSolucion.Range("B56").Value = "S052"
$wsSolucion.Range("A57").Value = "Pedido_19"
$wsSolucion.Range("B57").Value = "S016"
$wsSolucion.Range("A58").Value = "Pedido_26"
$wsSolucion.Range("B58").Value = "S075"
$wsSolucion.Range("A59").Value = "Pedido_68"
$wsSolucion.Range("B59").Value = "S039"
$wsSolucion.Range("A60").Value = "Pedido_24"
$wsSolucion.Range("B60").Value = "S056"
$wsSolucion.Range("A61").Value = "Pedido_6"
$wsSolucion.Range("B61").Value = "S036"
$wsSolucion.Range("A62").Value = "Pedido_15"
$wsSolucion.Range("B62").Value = "S017"
$wsSolucion.Range("A63").Value = "Pedido_9"
$wsSolucion.Range("B63").Value = "S057"
$wsSolucion.Range("A64").Value = "Pedido_70"
$wsSolucion.Range("B64").Value = "S079"
$wsSolucion.Range("A65").Value = "Pedido_7"
$wsSolucion.Range("B65").Value = "S040"
$wsSolucion.Range("A66").Value = "Pedido_20"
$wsSolucion.Range("B66").Value = "S021"
$wsSolucion.Range("A67").Value = "Pedido_56"
$wsSolucion.Range("B67").Value = "S076"
$wsSolucion.Range("A68").Value = "Pedido_54"
$wsSolucion.Range("B68").Value = "S061"
$wsSolucion.Range("A69").Value = "Pedido_55"
$wsSolucion.Range("B69").Value = "S058"
$wsSolucion.Range("A70").Value = "Pedido_27"
$wsSolucion.Range("B70").Value = "S080"
$wsSolucion.Range("A71").Value = "Pedido_31"
$wsSolucion.Range("B71").Value = "S018"
$wsSolucion.Range("A72").Value = "Pedido_2"
$wsSolucion.Range("B72").Value = "S062"
$wsSolucion.Range("A73").Value = "Pedido_41"
$wsSolucion.Range("B73").Value = "S022"
$wsSolucion.Range("A74").Value = "Pedido_53"
$wsSolucion.Range("B74").Value = "S019"
$wsSolucion.Range("A75").Value = "Pedido_76"
$wsSolucion.Range("B75").Value = "S059"
$wsSolucion.Range("A76").Value = "Pedido_60"
$wsSolucion.Range("B76").Value = "S023"
$wsSolucion.Range("A77").Value = "Pedido_72"
$wsSolucion.Range("B77").Value = "S063"
$wsSolucion.Range("A78").Value = "Pedido_44"
$wsSolucion.Range("B78").Value = "S060"
$wsSolucion.Range("A79").Value = "Pedido_14"
$wsSolucion.Range("B79").Value = "S020"
$wsSolucion.Range("A80").Value = "Pedido_11"
$wsSolucion.Range("B80").Value = "S064"
$wsSolucion.Range("A81").Value = "Pedido_38"
$wsSolucion.Range("B81").Value = "S024"

# --- Sheet "Metricas": update Tiempo (B2:B5) ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 693.3515515190833
$wsMetricas.Range("B3").Value = 448.8062493242512
$wsMetricas.Range("B4").Value = 693.987917612715
$wsMetricas.Range("B5").Value = 498.7823278192237
